$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged), update metric values
$ws.Range("B3").Value = 0.9908759461678647
$ws.Range("C3").Value = 0.9904333906452124
$ws.Range("D3").Value = 0.990433437944477

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9920399284300464
$ws.Range("C4").Value = 0.9920400393746225
$ws.Range("D4").Value = 0.9962316356274609

# Row 5: AdaBoostRegressor -> MLPRegressor, update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.6397366306640633
$ws.Range("C5").Value = 0.7455369856315704
$ws.Range("D5").Value = 0.6982310793695026
